$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

$ws.Range("H1").Value = "Save"

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
